$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference style (date column D) taken from an existing data row so the
# new date cells inherit the same number format (style index used for dates).
$dateFormat = $ws.Cells.Item(72, 4).NumberFormat

$newRows = @(
    @{
        A = 6
        B = "Mercado Mayorista Lo Valledor de Santiago"
        C = "Metropolitana"
        D = 44656
        E = 13
        F = "Fruta"
        G = 100104
        H = "Frutos de pepita"
        I = 100104003
        J = "Membrillo"
        K = "Champion"
        L = "Especial"
        M = 5
        N = 280000
        O = 280000
        P = 280000
        Q = "`$/bins (400 kilos)"
        R = "Región de O'Higgins"
        S = 700
        T = 400
    },
    @{
        A = 6
        B = "Mercado Mayorista Lo Valledor de Santiago"
        C = "Metropolitana"
        D = 44656
        E = 13
        F = "Fruta"
        G = 100104
        H = "Frutos de pepita"
        I = 100104003
        J = "Membrillo"
        K = "Champion"
        L = "Primera"
        M = 8
        N = 250000
        O = 250000
        P = 250000
        Q = "`$/bins (400 kilos)"
        R = "Región de O'Higgins"
        S = 625
        T = 400
    },
    @{
        A = 6
        B = "Mercado Mayorista Lo Valledor de Santiago"
        C = "Metropolitana"
        D = 44656
        E = 13
        F = "Fruta"
        G = 100104
        H = "Frutos de pepita"
        I = 100104003
        J = "Membrillo"
        K = "Champion"
        L = "Segunda"
        M = 12
        N = 220000
        O = 220000
        P = 220000
        Q = "`$/bins (400 kilos)"
        R = "Región de O'Higgins"
        S = 550
        T = 400
    }
)

$startRow = 73
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = $data.A
    $ws.Cells.Item($r, 2).Value = $data.B
    $ws.Cells.Item($r, 3).Value = $data.C

    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 4).NumberFormat = $dateFormat

    $ws.Cells.Item($r, 5).Value = $data.E
    $ws.Cells.Item($r, 6).Value = $data.F
    $ws.Cells.Item($r, 7).Value = $data.G
    $ws.Cells.Item($r, 8).Value = $data.H
    $ws.Cells.Item($r, 9).Value = $data.I
    $ws.Cells.Item($r, 10).Value = $data.J
    $ws.Cells.Item($r, 11).Value = $data.K
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
}
